# ---------------------------------------------------------------------------
# Applies the "Add files via upload" revision to ca2_s22.docx:
#   1. Rewrites the "(Due: ... )" italic line to "(Due Thursday 4/14/22,
#      11:30 PM)" and stamps an East-Asian font (Calibri) onto that line's
#      runs.
#   2. Moves the hidden "_GoBack" bookmark from the due-date line down to
#      wrap the "A uniform distributed load ... Table 1." sentence.
#   3. Relabels the load-span points from "A and B" to "B and C" in that
#      sentence (the "simply supported at points A and B" sentence stays
#      untouched).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ===========================================================================
# 1. Due-date line
# ===========================================================================

$dueRange = $d.Paragraphs(4).Range

# Stamp the East Asian font onto every run already in the line before the
# text is touched, so the untouched "(" run picks it up too.
$dueRange.Font.NameFarEast = "Calibri"

# "Due: " -> "Due "
$f1 = $d.Paragraphs(4).Range
$f1.Find.ClearFormatting()
$f1.Find.Execute("Due: ") | Out-Null
$f1.Text = "Due "
$f1.Font.Italic = $true
$f1.Font.NameFarEast = "Calibri"

# Insert "Thursday" right before the date
$f2 = $d.Paragraphs(4).Range
$f2.Find.ClearFormatting()
$f2.Find.Execute("4/14") | Out-Null
$f2.Collapse(1)
$f2.InsertBefore("Thursday")
$f2.Font.Italic = $true
$f2.Font.NameFarEast = "Calibri"

# "4/14/2022" -> " 4/14/22" (leading space, 2-digit year)
$f3 = $d.Paragraphs(4).Range
$f3.Find.ClearFormatting()
$f3.Find.Execute("4/14/2022") | Out-Null
$f3.Text = " 4/14/22"
$f3.Font.Italic = $true
$f3.Font.NameFarEast = "Calibri"

# " at 11:30 PM)" -> ", 11:30 PM)"
$f4 = $d.Paragraphs(4).Range
$f4.Find.ClearFormatting()
$f4.Find.Execute(" at 11:30 PM)") | Out-Null
$f4.Text = ", 11:30 PM)"
$f4.Font.Italic = $true
$f4.Font.NameFarEast = "Calibri"

# ===========================================================================
# 2. Move the "_GoBack" bookmark
# ===========================================================================

$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$sentenceStart = $d.Content
$sentenceStart.Find.ClearFormatting()
$sentenceStart.Find.Execute("A uniform distributed load") | Out-Null
$startPos = $sentenceStart.Start

$loadParagraph = $d.Paragraphs.Item(6).Range
$newBookmarkRange = $d.Range($startPos, $loadParagraph.End)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)

# ===========================================================================
# 3. "points A and B" -> "points B and C" (load-span sentence only)
# ===========================================================================

$betweenPoints = $d.Content
$betweenPoints.Find.ClearFormatting()
$betweenPoints.Find.Execute("is applied between points ") | Out-Null
$afterPoints = $betweenPoints.End

$letterA = $d.Range($afterPoints, $afterPoints + 1)
$letterA.Text = "B"

$andRange = $d.Range($afterPoints, $afterPoints + 20)
$andRange.Find.ClearFormatting()
$andRange.Find.Execute(" and ") | Out-Null
$posAfterAnd = $andRange.End

$letterB = $d.Range($posAfterAnd, $posAfterAnd + 1)
$letterB.Text = "C"
